$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("ZZ1")

$ws.Range("D2").Value = '42.774.49'
$ws.Range("E2").Value = '  -2.11%  '
$ws.Range("D3").Value = '2.333.40'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  +0.04%  '
$helper.NumberFormat = "@"
$helper.Value = '306.41'
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -2.10%  '
$helper.NumberFormat = "@"
$helper.Value = '99.92'
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  -3.42%  '
$ws.Range("E7").Value = '  -5.48%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -4.71%  '
$helper.NumberFormat = "@"
$helper.Value = '34.94'
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -3.60%  '
$helper.NumberFormat = "@"
$helper.Value = '52.10'
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("E12").Value = '  -2.52%  '
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("E14").Value = '  -3.68%  '
$helper.NumberFormat = "@"
$helper.Value = '15.77'
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  +4.33%  '
$ws.Range("D16").Value = '2.369.01'
$ws.Range("E16").Value = '  +2.01%  '
$helper.NumberFormat = "@"
$helper.Value = '0.797'
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  -2.58%  '
$ws.Range("D18").Value = '42.714.14'
$ws.Range("E18").Value = '  -1.98%  '
$helper.NumberFormat = "@"
$helper.Value = '6.23'
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  +0.81%  '
$ws.Range("D20").Value = '0.0₃0907'
$ws.Range("E20").Value = '  -3.11%  '
$helper.NumberFormat = "@"
$helper.Value = '11.64'
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -7.60%  '
$helper.NumberFormat = "@"
$helper.Value = '67.75'
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  -1.16%  '
$helper.NumberFormat = "@"
$helper.Value = '236.72'
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  -2.68%  '
$ws.Range("E24").Value = '  -3.15%  '
$helper.NumberFormat = "@"
$helper.Value = '2.57'
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  -2.99%  '
$ws.Range("E26").Value = '  +0.02%  '
$helper.NumberFormat = "@"
$helper.Value = '25.10'
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +0.99%  '
$helper.NumberFormat = "@"
$helper.Value = '2.32'
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$helper.NumberFormat = "@"
$helper.Value = '34.79'
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  -8.09%  '
$helper.NumberFormat = "@"
$helper.Value = '9.35'
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  -3.66%  '
$helper.NumberFormat = "@"
$helper.Value = '159.32'
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -5.11%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("E33").Value = '  -4.48%  '
$ws.Range("B34").Value = 'Celestia'
$ws.Range("C34").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$helper.NumberFormat = "@"
$helper.Value = '17.37'
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  -2.57%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$helper.NumberFormat = "@"
$helper.Value = '2.45'
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -3.30%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$helper.NumberFormat = "@"
$helper.Value = '0.0727'
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -3.11%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$helper.NumberFormat = "@"
$helper.Value = '4.58'
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  +5.13%  '
$helper.NumberFormat = "@"
$helper.Value = '2.96'
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  -5.96%  '
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("E40").Value = '  -4.08%  '
$helper.NumberFormat = "@"
$helper.Value = '0.113'
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -3.45%  '
$helper.NumberFormat = "@"
$helper.Value = '2.34'
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  +1.35%  '
$ws.Range("D43").Value = '2.023.72'
$ws.Range("E43").Value = '  +1.95%  '
$helper.NumberFormat = "@"
$helper.Value = '0.0285'
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$helper.NumberFormat = "@"
$helper.Value = '18.69'
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  -6.63%  '
$helper.NumberFormat = "@"
$helper.Value = '10.24'
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  +3.11%  '
$ws.Range("E47").Value = '  -3.69%  '
$helper.NumberFormat = "@"
$helper.Value = '56.17'
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("E49").Value = '  -1.14%  '
$ws.Range("D50").Value = '2.559.10'
$ws.Range("E50").Value = '  +0.38%  '
$ws.Range("E51").Value = '  +1.59%  '

$helper.Clear()
